$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1309.5714
$ws.Range("I137").Value = 1089.1428
$ws.Range("J137").Value = 1970.8572
$ws.Range("K137").Value = 3267.4284
$ws.Range("L137").Value = 5912.571599999999
$ws.Range("M137").Value = -717.4284000000002
$ws.Range("N137").Value = -11012.5716
$ws.Range("H141").Value = 5054.8335
$ws.Range("I141").Value = 4665.8
$ws.Range("K141").Value = 13997.4
$ws.Range("M141").Value = -8817.400000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 26316896
$ws.Range("I2").Value = 34483668
$ws.Range("K2").Value = 34483668
$ws.Range("M2").Value = -34483555
$ws.Range("H61").Value = 2943.305
$ws.Range("I61").Value = 1947.5333
$ws.Range("J61").Value = 5659.0454
$ws.Range("K61").Value = 1947.5333
$ws.Range("L61").Value = 5659.0454
$ws.Range("M61").Value = -1735.5333
$ws.Range("N61").Value = -6083.0454
$ws.Range("H63").Value = 1991.8889
$ws.Range("I63").Value = 1843.7142
$ws.Range("K63").Value = 1843.7142
$ws.Range("M63").Value = -1157.7142
$ws.Range("H66").Value = 1991.8889
$ws.Range("I66").Value = 1843.7142
$ws.Range("K66").Value = 9218.571
$ws.Range("M66").Value = -5786.571
$ws.Range("H74").Value = 8959.462
$ws.Range("J74").Value = 24665.875
$ws.Range("L74").Value = 24665.875
$ws.Range("N74").Value = -26413.875
$ws.Range("H77").Value = 8959.462
$ws.Range("J77").Value = 24665.875
$ws.Range("L77").Value = 123329.375
$ws.Range("N77").Value = -132065.375
$ws.Range("H110").Value = 6372.1665
$ws.Range("I110").Value = 6425.381
$ws.Range("K110").Value = 6425.381
$ws.Range("M110").Value = -4380.381
$ws.Range("H116").Value = 26316896
$ws.Range("I116").Value = 34483668
$ws.Range("K116").Value = 34483668
$ws.Range("M116").Value = -34481374
$ws.Range("H132").Value = 2366.6135
$ws.Range("I132").Value = 2131.5642
$ws.Range("K132").Value = 6394.692599999999
$ws.Range("M132").Value = -3864.692599999999
$ws.Range("H136").Value = 2943.305
$ws.Range("I136").Value = 1947.5333
$ws.Range("J136").Value = 5659.0454
$ws.Range("K136").Value = 5842.5999
$ws.Range("L136").Value = 16977.1362
$ws.Range("M136").Value = -3292.5999
$ws.Range("N136").Value = -22077.1362
$ws.Range("H139").Value = 95499.8
$ws.Range("J139").Value = 95499.8
$ws.Range("L139").Value = 95499.8
$ws.Range("N139").Value = -105779.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 26316896
$ws.Range("I3").Value = 34483668
$ws.Range("K3").Value = 34483668
$ws.Range("M3").Value = -34483554
$ws.Range("H11").Value = 129.66667
$ws.Range("I11").Value = 94.5
$ws.Range("J11").Value = 200
$ws.Range("K11").Value = 94.5
$ws.Range("L11").Value = 200
$ws.Range("M11").Value = 45.5
$ws.Range("N11").Value = -480
$ws.Range("H94").Value = 860.9
$ws.Range("I94").Value = 516.46155
$ws.Range("K94").Value = 516.46155
$ws.Range("M94").Value = -65.46154999999999
$ws.Range("H105").Value = 1193.6666
$ws.Range("I105").Value = 1230.8125
$ws.Range("K105").Value = 1230.8125
$ws.Range("M105").Value = 516.1875
$ws.Range("H107").Value = 1311.5714
$ws.Range("I107").Value = 1196.8334
$ws.Range("K107").Value = 1196.8334
$ws.Range("M107").Value = 723.1666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2541.75
$ws.Range("I58").Value = 2076.6
$ws.Range("J58").Value = 3937.2
$ws.Range("K58").Value = 2076.6
$ws.Range("L58").Value = 3937.2
$ws.Range("M58").Value = -1873.6
$ws.Range("N58").Value = -4343.2
$ws.Range("H107").Value = 622.0417
$ws.Range("I107").Value = 585.8125
$ws.Range("J107").Value = 694.5
$ws.Range("K107").Value = 585.8125
$ws.Range("L107").Value = 694.5
$ws.Range("M107").Value = 1334.1875
$ws.Range("N107").Value = -4534.5
$ws.Range("H136").Value = 2541.75
$ws.Range("I136").Value = 2076.6
$ws.Range("J136").Value = 3937.2
$ws.Range("K136").Value = 6229.799999999999
$ws.Range("L136").Value = 11811.6
$ws.Range("M136").Value = -3679.799999999999
$ws.Range("N136").Value = -16911.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 1198.8
$ws.Range("J107").Value = 1299
$ws.Range("L107").Value = 3897
$ws.Range("N107").Value = -7737
$ws.Range("H133").Value = 23557
$ws.Range("J133").Value = 13500
$ws.Range("L133").Value = 40500
$ws.Range("N133").Value = -50620
$ws.Range("H134").Value = 7606.467
$ws.Range("I134").Value = 5508.0835
$ws.Range("J134").Value = 16000
$ws.Range("K134").Value = 16524.2505
$ws.Range("L134").Value = 48000
$ws.Range("M134").Value = -11454.2505
$ws.Range("N134").Value = -58140
$ws.Range("H139").Value = 2130.25
$ws.Range("I139").Value = 1851.4445
$ws.Range("K139").Value = 5554.333500000001
$ws.Range("M139").Value = -414.3335000000006

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7418.6665
$ws.Range("I70").Value = 7856.8237
$ws.Range("J70").Value = 6845.6924
$ws.Range("K70").Value = 7856.8237
$ws.Range("L70").Value = 6845.6924
$ws.Range("M70").Value = -7586.8237
$ws.Range("N70").Value = -7385.6924
$ws.Range("H73").Value = 7418.6665
$ws.Range("I73").Value = 7856.8237
$ws.Range("J73").Value = 6845.6924
$ws.Range("K73").Value = 7856.8237
$ws.Range("L73").Value = 6845.6924
$ws.Range("M73").Value = -6920.8237
$ws.Range("N73").Value = -8717.6924
$ws.Range("H126").Value = 21167.334
$ws.Range("I126").Value = 29938.625
$ws.Range("K126").Value = 89815.875
$ws.Range("M126").Value = -87345.875
$ws.Range("H132").Value = 3091.9473
$ws.Range("I132").Value = 2831.6875
$ws.Range("J132").Value = 4480
$ws.Range("K132").Value = 8495.0625
$ws.Range("L132").Value = 13440
$ws.Range("M132").Value = -5965.0625
$ws.Range("N132").Value = -18500

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 105650
$ws.Range("I61").Value = 116611.11
$ws.Range("K61").Value = 116611.11
$ws.Range("M61").Value = -116409.11
$ws.Range("H113").Value = 105650
$ws.Range("I113").Value = 116611.11
$ws.Range("K113").Value = 116611.11
$ws.Range("M113").Value = -114441.11
$ws.Range("H132").Value = 6240.6
$ws.Range("I132").Value = 4404.3335
$ws.Range("K132").Value = 13213.0005
$ws.Range("M132").Value = -10683.0005
$ws.Range("H136").Value = 4222.1665
$ws.Range("I136").Value = 3806.9473
$ws.Range("K136").Value = 11420.8419
$ws.Range("M136").Value = -8870.841899999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1529.1
$ws.Range("I132").Value = 1404.5555
$ws.Range("K132").Value = 4213.666499999999
$ws.Range("M132").Value = -1683.666499999999
$ws.Range("H136").Value = 2032.2903
$ws.Range("I136").Value = 1777.6296
$ws.Range("K136").Value = 5332.8888
$ws.Range("M136").Value = -2782.8888
